# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Kraken Profits workbook. Updates cached price/profit values (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 209.57143
$ws.Range("I9").Value = 166.75
$ws.Range("J9").Value = 266.66666
$ws.Range("K9").Value = 166.75
$ws.Range("L9").Value = 266.66666
$ws.Range("M9").Value = 2.25
$ws.Range("N9").Value = -604.66666
$ws.Range("H92").Value = 679.8
$ws.Range("I92").Value = 466.33334
$ws.Range("K92").Value = 466.33334
$ws.Range("M92").Value = 781.66666
$ws.Range("H129").Value = 1260.2858
$ws.Range("I129").Value = 1260.2858
$ws.Range("K129").Value = 3780.8574
$ws.Range("M129").Value = 1219.1426
$ws.Range("H132").Value = 5615.364
$ws.Range("I132").Value = 3602.4285
$ws.Range("J132").Value = 9138
$ws.Range("K132").Value = 10807.2855
$ws.Range("L132").Value = 27414
$ws.Range("M132").Value = -8277.2855
$ws.Range("N132").Value = -32474
$ws.Range("H135").Value = 1416.3846
$ws.Range("I135").Value = 1101.6666
$ws.Range("K135").Value = 9914.999400000001
$ws.Range("M135").Value = -7379.999400000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H38").Value = 2000
$ws.Range("I38").Value = 2000
$ws.Range("K38").Value = 2000
$ws.Range("M38").Value = -1533
$ws.Range("H63").Value = 3155.8
$ws.Range("I63").Value = 3155.8
$ws.Range("K63").Value = 3155.8
$ws.Range("M63").Value = -2469.8
$ws.Range("H66").Value = 3155.8
$ws.Range("I66").Value = 3155.8
$ws.Range("K66").Value = 15779
$ws.Range("M66").Value = -12347
$ws.Range("H88").Value = 1200
$ws.Range("J88").Value = 1500
$ws.Range("L88").Value = 1500
$ws.Range("N88").Value = -2312
$ws.Range("H91").Value = 1200
$ws.Range("J91").Value = 1500
$ws.Range("L91").Value = 1500
$ws.Range("N91").Value = -4308
$ws.Range("H102").Value = 5999.5
$ws.Range("I102").Value = 5999.5
$ws.Range("K102").Value = 5999.5
$ws.Range("M102").Value = -4377.5
$ws.Range("H132").Value = 4728.4116
$ws.Range("I132").Value = 4930.1875
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 14790.5625
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -12260.5625
$ws.Range("N132").Value = -9560

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 350
$ws.Range("I20").Value = 350
$ws.Range("K20").Value = 350
$ws.Range("M20").Value = -103
$ws.Range("H33").Value = 22682
$ws.Range("I33").Value = 22011
$ws.Range("J33").Value = 24024
$ws.Range("K33").Value = 22011
$ws.Range("L33").Value = 24024
$ws.Range("M33").Value = -21675
$ws.Range("N33").Value = -24696
$ws.Range("H94").Value = 1531.6666
$ws.Range("I94").Value = 1297.5
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1297.5
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -846.5
$ws.Range("N94").Value = -2902
$ws.Range("H105").Value = 2601.4285
$ws.Range("I105").Value = 2302.5
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2302.5
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -555.5
$ws.Range("N105").Value = -6494

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5632
$ws.Range("H35").Value = 1318.5
$ws.Range("I35").Value = 1318.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1318.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1024.5
$ws.Range("N35").ClearContents()
$ws.Range("H38").Value = 3047
$ws.Range("I38").Value = 3047
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3047
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2670
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 3047
$ws.Range("I46").Value = 3047
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3047
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2836
$ws.Range("N46").ClearContents()
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21240
$ws.Range("H107").Value = 339.7647
$ws.Range("I107").Value = 277.16666
$ws.Range("J107").Value = 490
$ws.Range("K107").Value = 277.16666
$ws.Range("L107").Value = 490
$ws.Range("M107").Value = 1642.83334
$ws.Range("N107").Value = -4330
$ws.Range("H134").Value = 1292.6666
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2688.7144
$ws.Range("J131").Value = 4524.75
$ws.Range("L131").Value = 13574.25
$ws.Range("N131").Value = -23654.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28.666666
$ws.Range("I2").Value = 25.052631
$ws.Range("K2").Value = 25.052631
$ws.Range("M2").Value = 87.94736899999999
$ws.Range("H70").Value = 4999.5
$ws.Range("J70").Value = 4999.5
$ws.Range("L70").Value = 4999.5
$ws.Range("N70").Value = -5539.5
$ws.Range("H73").Value = 4999.5
$ws.Range("J73").Value = 4999.5
$ws.Range("L73").Value = 4999.5
$ws.Range("N73").Value = -6871.5
$ws.Range("H97").Value = 4497.5
$ws.Range("J97").Value = 7000
$ws.Range("L97").Value = 7000
$ws.Range("N97").Value = -7992
$ws.Range("H102").Value = 3656.7144
$ws.Range("I102").Value = 3656.7144
$ws.Range("K102").Value = 3656.7144
$ws.Range("M102").Value = -2034.7144
$ws.Range("H122").Value = 7608.2856
$ws.Range("I122").Value = 2995.8333
$ws.Range("K122").Value = 8987.499899999999
$ws.Range("M122").Value = -6537.499899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2621.3333
$ws.Range("I93").Value = 2574
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2574
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1326
$ws.Range("N93").Value = -5496

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6442.3076
$ws.Range("I100").Value = 7824.4
$ws.Range("J100").Value = 1835.3334
$ws.Range("K100").Value = 15648.8
$ws.Range("L100").Value = 3670.6668
$ws.Range("M100").Value = -15107.8
$ws.Range("N100").Value = -4752.6668
$ws.Range("H113").Value = 477.91666
$ws.Range("I113").Value = 459.5
$ws.Range("K113").Value = 1378.5
$ws.Range("M113").Value = 791.5
$ws.Range("H136").Value = 2004
$ws.Range("I136").Value = 1759.45
$ws.Range("K136").Value = 5278.35
$ws.Range("M136").Value = -2728.35
